$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 is a duplicate of existing row 3.
# New row 6 is a duplicate of existing row 2.
# New row 7 is a duplicate of existing row 4.
$sourceForNewRow = @{ 5 = 3; 6 = 2; 7 = 4 }

foreach ($newRow in 5..7) {
    $srcRow = $sourceForNewRow[$newRow]
    for ($col = 1; $col -le 8; $col++) {
        $srcCell = $ws.Cells.Item($srcRow, $col)
        $dstCell = $ws.Cells.Item($newRow, $col)
        $text = $srcCell.Text

        # Numeric-looking strings (e.g. "0303") must stay text, otherwise
        # Excel coerces them into numbers and the leading zero is lost -
        # force the destination cell to Text format before writing.
        if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
            $dstCell.NumberFormat = "@"
        }

        $dstCell.Value = $text
    }
}
